# DataSource - Emision Motor_Answer.xlsx
# "se sube para emitir las pólizas de movilidad en QA"
#
# Updates row 2 (Fecha/Cuenta) and row 3 (Ambiente/URL/Fecha/Campania)
# of the data source sheet with the new QA environment values, and
# moves the on-screen selection over to column W.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: change a cell's value while keeping its current number format /
# font / style (Excel tends to reset quotePrefix / custom numFmt driven
# styles back to the default style whenever .Value is assigned directly).
function Set-CellValuePreserveFormat {
    param(
        $Range,
        $Value
    )

    $holding = $ws.Range("ZZ1000")
    $Range.Copy() | Out-Null
    $holding.PasteSpecial(-4122) | Out-Null   # xlPasteFormats
    $Range.Value = $Value
    $holding.Copy() | Out-Null
    $Range.PasteSpecial(-4122) | Out-Null     # xlPasteFormats
    $holding.Clear() | Out-Null
    $excel.CutCopyMode = $false
}

# ---------------------------------------------------------------------
# Row 2: cuenta / fecha de inicio de la poliza SPA001
# ---------------------------------------------------------------------
Set-CellValuePreserveFormat $ws.Range("F2") 1728150905
Set-CellValuePreserveFormat $ws.Range("J2") "26/04/2021"

# ---------------------------------------------------------------------
# Row 3: ambiente / URL / fecha / campania RGM006 (antes RGA006)
# ---------------------------------------------------------------------
$ws.Range("B3").Value = "i-preproducciongestion.segurossura.com.ar"
$ws.Range("C3").Value = "https://i-preproducciongestion.segurossura.com.ar/pc/PolicyCenter.do"
Set-CellValuePreserveFormat $ws.Range("J3") "03/05/2021"
$ws.Range("U3").Value = "RGM006"
$ws.Range("V3").Value = "ABC12RGM006"
$ws.Range("W3").Value = "ZAZ123RGM006"

# New cell added below the campania (blank placeholder with a single space)
$ws.Range("V4").Value = " "

# ---------------------------------------------------------------------
# Window/view state: scroll right and select W4, matching the author's
# last recorded selection in the workbook.
# ---------------------------------------------------------------------
$ws.Range("M1").Select() | Out-Null
$excel.ActiveWindow.ScrollColumn = 13
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("W4").Select() | Out-Null
